$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Sheet1's selection becomes A:C (full columns A-C), and it loses the active-tab flag
# once Sheet2 is added after it and becomes active.
$ws1.Range("A:C").Select()

# Add "Sheet2" right after "Sheet1"
$ws2 = $wb.Worksheets.Add($null, $ws1)

# Header row
$ws2.Range("A1").Value = "ID"
$ws2.Range("B1").Value = "age"
$ws2.Range("C1").Value = "sex"

# Data rows (sub_1..sub_50, age, sex)
$data = @(
    @("sub_1", 54.3893182249734, "Male"),
    @("sub_2", 62.5324198490171, "Male"),
    @("sub_3", 74.6096683642242, "Female"),
    @("sub_4", 93.0878075945433, "Male"),
    @("sub_5", 93.4199840839674, "Female"),
    @("sub_6", 57.0925886754896, "Female"),
    @("sub_7", 93.6766751792277, "Female"),
    @("sub_8", 93.0725126709325, "Male"),
    @("sub_9", 71.8419041925278, "Male"),
    @("sub_10", 86.012621099996, "Female"),
    @("sub_11", 56.3848852382246, "Male"),
    @("sub_12", 68.9792577181823, "Female"),
    @("sub_13", 91.208098633508, "Male"),
    @("sub_14", 85.6493298301799, "Female"),
    @("sub_15", 93.1771591876806, "Male"),
    @("sub_16", 79.5083314620464, "Female"),
    @("sub_17", 51.6070255358385, "Female"),
    @("sub_18", 88.2108187640949, "Female"),
    @("sub_19", 92.0296961490897, "Male"),
    @("sub_20", 80.5430819685998, "Male"),
    @("sub_21", 84.0983058760249, "Female"),
    @("sub_22", 83.4409610656212, "Female"),
    @("sub_23", 67.6502158790375, "Male"),
    @("sub_24", 79.49650505799, "Female"),
    @("sub_25", 57.7034009515202, "Female"),
    @("sub_26", 81.7720739608823, "Male"),
    @("sub_27", 51.4324780869839, "Female"),
    @("sub_28", 62.46153432324, "Male"),
    @("sub_29", 52.0777125784019, "Male"),
    @("sub_30", 54.3709301556131, "Female"),
    @("sub_31", 87.0556022747281, "Female"),
    @("sub_32", 81.2672880339117, "Male"),
    @("sub_33", 64.2694766027387, "Female"),
    @("sub_34", 92.7599921977259, "Male"),
    @("sub_35", 51.5500736226308, "Male"),
    @("sub_36", 69.7434961845379, "Male"),
    @("sub_37", 67.1701305691853, "Male"),
    @("sub_38", 84.4482554667051, "Male"),
    @("sub_39", 85.7839955511678, "Male"),
    @("sub_40", 58.409267204947, "Male"),
    @("sub_41", 72.0393978104704, "Male"),
    @("sub_42", 70.0513790319904, "Male"),
    @("sub_43", 79.0840854550069, "Female"),
    @("sub_44", 81.9214173886132, "Male"),
    @("sub_45", 83.9609006892062, "Female"),
    @("sub_46", 62.421128464936, "Male"),
    @("sub_47", 80.5866204584153, "Male"),
    @("sub_48", 79.4794101788228, "Male"),
    @("sub_49", 57.3175280837583, "Male"),
    @("sub_50", 55.3548956701269, "Female")
)

$r = 2
foreach ($row in $data) {
    $ws2.Cells.Item($r, 1).Value = $row[0]
    $ws2.Cells.Item($r, 2).Value = $row[1]
    $ws2.Cells.Item($r, 3).Value = $row[2]
    $r++
}

# --- Formatting ---

# Header row: bold, centered horizontally and vertically
$header = $ws2.Range("A1:C1")
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4108

# Rows 2-6, columns A and B: centered horizontally and vertically
$ws2.Range("A2:B6").HorizontalAlignment = -4108
$ws2.Range("A2:B6").VerticalAlignment = -4108

# Rows 7-51, column A: centered horizontally, bottom vertically
$ws2.Range("A7:A51").HorizontalAlignment = -4108
$ws2.Range("A7:A51").VerticalAlignment = -4107

# Column C, rows 2-51 ("sex"): centered horizontally, bottom vertically, text format
$sexCol = $ws2.Range("C2:C51")
$sexCol.HorizontalAlignment = -4108
$sexCol.VerticalAlignment = -4107
$sexCol.NumberFormat = "@"

# Select A:C on the new sheet (matches the saved selection state) and make it active
$ws2.Range("A:C").Select()

Write-Output "Sheet2 added with $($data.Count) data rows"
